$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table (Tipo Doc / N Doc / Nombre / Periodo Mora / Valor Mora / Salario Basico) occupies rows 16-18.
# The underlying database was corrected/reordered: row 16 now holds JORGE LUIS COVO MARTINEZ's data
# (with a corrected Valor Mora), row 17 now holds PAOLA ANDREA LUENGAS TORRES's data, and row 18 now
# holds JAIR ALI GONZALEZ SANTIAGO's data.

# Row 16 -> JORGE LUIS COVO MARTINEZ
$ws.Range("C16").Value = "1129517708"
$ws.Range("D16").Value = "JORGE LUIS COVO MARTINEZ"
$ws.Range("E16").Value = "1908"
$ws.Range("F16").Value = 66250
$ws.Range("G16").Value = 1656232

# Row 17 -> PAOLA ANDREA LUENGAS TORRES
$ws.Range("C17").Value = "1102870612"
$ws.Range("D17").Value = "PAOLA ANDREA LUENGAS TORRES"
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 36000
$ws.Range("G17").Value = 955790

# Row 18 -> JAIR ALI GONZALEZ SANTIAGO
$ws.Range("C18").Value = "72162724"
$ws.Range("D18").Value = "JAIR ALI GONZALEZ SANTIAGO"
$ws.Range("E18").Value = "1908"
$ws.Range("F18").Value = 66250
$ws.Range("G18").Value = 1656232
